# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" worksheets to reflect the latest scrape.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 8371
    "F3"  = 7893
    "F8"  = 132
    "F14" = 1882
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
